# Updates cryptos list values (Price and Volume(1h) columns) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) cells whose new text would otherwise be auto-parsed as a number ---
# Force these to remain plain text (matching the original inlineStr cell type) by
# temporarily applying a text number format, then restoring the default "Normal" style
# so no stray style reference is left behind.
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '212.22'
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.491'
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '46.42'
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '24.37'
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.0882'
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '62.28'
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '228.27'
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '7.39'
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '3.91'
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '9.14'
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.02'
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '151.27'
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '15.00'
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '6.46'
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.104'
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '3.21'
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '3.14'
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.532'
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.86'
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.979'
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '86.08'
$cell.Style = "Normal"

# --- Price (column D) cells whose new text is not number-like; plain assignment is safe ---
$ws.Range("D2").Value = '28.575.14'
$ws.Range("D3").Value = '1.575.32'
$ws.Range("D13").Value = '1.802.16'
$ws.Range("D14").Value = '1.571.15'
$ws.Range("D17").Value = '28.569.01'
$ws.Range("D21").Value = '0.0₃0694'
$ws.Range("D35").Value = '1.394.84'
$ws.Range("D48").Value = '1.713.89'

# --- Volume(1h) (column E) cells; values contain "%" and spaces so they stay text ---
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  +6.60%  '
$ws.Range("E9").Value = '  +4.38%  '
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("E11").Value = '  -1.29%  '
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("E18").Value = '  -2.22%  '
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("E21").Value = '  -1.90%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  -5.26%  '
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("E25").Value = '  +3.62%  '
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("E28").Value = '  -1.96%  '
$ws.Range("E29").Value = '  -2.10%  '
$ws.Range("E30").Value = '  +0.28%  '
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("E32").Value = '  -2.10%  '
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("E34").Value = '  +0.34%  '
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("E36").Value = '  -2.86%  '
$ws.Range("E37").Value = '  -2.24%  '
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("E39").Value = '  +4.95%  '
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("E45").Value = '  +2.59%  '
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  -4.06%  '
$ws.Range("E51").Value = '  -1.28%  '

